# Apply the commit's changes:
#  1. Insert a new "Player Info" worksheet as the first sheet, with
#     ID / NAME / BATTING_HAND / BOWL_STYLE columns for player 4683
#     (Mohammad Abbas).
#  2. Rename the MATCH_CARD_LINK column to MATCH_CODE on both the
#     "ODI Batting" and "ODI Bowling" sheets, replacing the full
#     scorecard URL values with the bare numeric match code (stored as
#     text, same as the rest of that column).
#
# NOTE: Worksheets.Add() inserts the new sheet in front of every other
# sheet, which shifts sheet *positions* - so any worksheet reference
# obtained before the Add() call can silently end up pointing at the
# wrong sheet afterwards. To stay safe we do the structural change
# (add + rename) FIRST, and only look up "ODI Batting"/"ODI Bowling"
# by name AFTER that, before doing any value edits.

$wb = $excel.ActiveWorkbook

# --- 1. New "Player Info" sheet, inserted before all existing sheets ---
$info = $wb.Worksheets.Add()
$info.Name = "Player Info"

# Look up the other two sheets by name now that the sheet list is final.
$batting = $wb.Worksheets.Item("ODI Batting")
$bowling = $wb.Worksheets.Item("ODI Bowling")

# Copy the header formatting (bold, bordered, centered) from an existing
# header cell so the new headers match the rest of the workbook, then
# overwrite the copied values with the real header text.
$batting.Range("A1").Copy($info.Range("A1"))
$batting.Range("A1").Copy($info.Range("B1"))
$batting.Range("A1").Copy($info.Range("C1"))
$batting.Range("A1").Copy($info.Range("D1"))

$info.Range("A1").Value = "ID"
$info.Range("B1").Value = "NAME"
$info.Range("C1").Value = "BATTING_HAND"
$info.Range("D1").Value = "BOWL_STYLE"

# Data row - player id is textual (matches the workbook's existing
# convention of storing ids/codes as text), so force text with a leading
# apostrophe rather than letting it be parsed as a number, then drop the
# resulting "quote prefix" style tweak so the cell stays on the default,
# unstyled format (matching every other plain data cell in this file).
$info.Range("A2").Value = "'4683"
$info.Range("A2").Style = "Normal"
$info.Range("B2").Value = "Mohammad Abbas"
$info.Range("C2").Value = "Right Handed"
$info.Range("D2").Value = "Right Arm Fast"

# --- 2. ODI Batting: MATCH_CARD_LINK -> MATCH_CODE ---
$batting.Range("D1").Value = "MATCH_CODE"
$batting.Range("D2").Value = "'4273"
$batting.Range("D2").Style = "Normal"
$batting.Range("D3").Value = "'4274"
$batting.Range("D3").Style = "Normal"
$batting.Range("D4").Value = "'4277"
$batting.Range("D4").Style = "Normal"

# --- 3. ODI Bowling: MATCH_CARD_LINK -> MATCH_CODE ---
$bowling.Range("B1").Value = "MATCH_CODE"
$bowling.Range("B2").Value = "'4273"
$bowling.Range("B2").Style = "Normal"
$bowling.Range("B3").Value = "'4274"
$bowling.Range("B3").Style = "Normal"
$bowling.Range("B4").Value = "'4277"
$bowling.Range("B4").Style = "Normal"

Write-Output "Player Info sheet added; MATCH_CARD_LINK columns converted to MATCH_CODE."
